# Handback status refresh: replace the two e2e test file UUIDs and refresh
# the "Generate Report for Handback" timestamps on all three sheets.
#
#   2b2094ba-a691-4556-8780-668ee1bd88b3  ->  be3fe9fd-d7f6-43a2-b38c-22c1a60f6fdd
#   6b220dc0-90e9-4acd-add0-3f4f76d7f157  ->  ffff793dbab2-2335-458b-b9da-9cc16614de8d
#
# Hyperlink targets (the git-blob URLs in the external rels) are left as-is;
# only the cell text / hyperlink display text is refreshed, matching the
# commit's xlsx OOXML diff.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "2b2094ba-a691-4556-8780-668ee1bd88b3"
$oldUuid2 = "6b220dc0-90e9-4acd-add0-3f4f76d7f157"
$newUuid1 = "be3fe9fd-d7f6-43a2-b38c-22c1a60f6fdd"
$newUuid2 = "ffff793dbab2-2335-458b-b9da-9cc16614de8d"

$addr1 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/2b05b1afd3129d27a4020040637d8bf3e70d10a4/e2e/$oldUuid1.md"
$addr2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/2b05b1afd3129d27a4020040637d8bf3e70d10a4/e2e/$oldUuid2.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("A3").Value = "$newUuid2.md"
$wsOverview.Range("G2").Value = "2016-08-14 01:29:43"
$wsOverview.Range("G3").Value = "2016-08-14 01:29:43"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $addr1, "", "", "e2e\$newUuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $addr2, "", "", "e2e\$newUuid2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhAddrA2 = $addr1
$zhAddrI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c2219b9faa87ece3079b184abae948da8ce3b06c/e2e/$oldUuid1.md"
$zhAddrA3 = $addr2
$zhAddrI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c2219b9faa87ece3079b184abae948da8ce3b06c/e2e/$oldUuid2.md"

$wsZhCn.Range("A2").Value = "$newUuid1.md"
$wsZhCn.Range("G2").Value = "$newUuid1.09cfd9f958ddaf58f738cb8355c698180b020c9f.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-14 01:29:35"
$wsZhCn.Range("I2").Value = "$newUuid1.md"
$wsZhCn.Range("J2").Value = "$newUuid1.09cfd9f958ddaf58f738cb8355c698180b020c9f.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-14 01:30:11"

$wsZhCn.Range("A3").Value = "$newUuid2.md"
$wsZhCn.Range("G3").Value = "$newUuid1.09cfd9f958ddaf58f738cb8355c698180b020c9f.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-14 01:29:35"
$wsZhCn.Range("I3").Value = "$newUuid2.md"
$wsZhCn.Range("J3").Value = "$newUuid1.09cfd9f958ddaf58f738cb8355c698180b020c9f.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-14 01:30:11"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhAddrA2, "", "", "$newUuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhAddrI2, "", "", "$newUuid1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $zhAddrA3, "", "", "$newUuid2.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhAddrI3, "", "", "$newUuid2.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deAddrA2 = $addr1
$deAddrI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/769d15bf13de3cd6c283cfe6117af83fcbc59ba5/e2e/$oldUuid1.md"
$deAddrA3 = $addr2
$deAddrI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/769d15bf13de3cd6c283cfe6117af83fcbc59ba5/e2e/$oldUuid2.md"

$wsDeDe.Range("A2").Value = "$newUuid1.md"
$wsDeDe.Range("G2").Value = "$newUuid1.09cfd9f958ddaf58f738cb8355c698180b020c9f.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-14 01:29:43"
$wsDeDe.Range("I2").Value = "$newUuid1.md"
$wsDeDe.Range("J2").Value = "$newUuid1.09cfd9f958ddaf58f738cb8355c698180b020c9f.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-14 01:30:22"

$wsDeDe.Range("A3").Value = "$newUuid2.md"
$wsDeDe.Range("G3").Value = "$newUuid1.09cfd9f958ddaf58f738cb8355c698180b020c9f.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-14 01:29:43"
$wsDeDe.Range("I3").Value = "$newUuid2.md"
$wsDeDe.Range("J3").Value = "$newUuid1.09cfd9f958ddaf58f738cb8355c698180b020c9f.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-14 01:30:22"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deAddrA2, "", "", "$newUuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $deAddrI2, "", "", "$newUuid1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $deAddrA3, "", "", "$newUuid2.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $deAddrI3, "", "", "$newUuid2.md")
